$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New simulation rows appended to the sheet (rows 66-80)
$newRows = @(
  @{ Row = 66; Cells = @("Elf", "[25", " -3]", 23, "S", "H", "H", "H", "S", "P", "S", "P", "S", "P", "S", "P", "S", "P", "S", "P", "S", "P", "S", "P", "S", "P", "S") }
  @{ Row = 67; Cells = @("Elf", "[25", " -3]", 13, "A", "H", "A", "H", "A", "P", "S", "P", "S", "P", "S", "P", "S") }
  @{ Row = 68; Cells = @("Elf", "[22", " -3]", 17, "A", "H", "A", "H", "A", "P", "A", "P", "S", "P", "S", "P", "S", "P", "S", "P", "S") }
  @{ Row = 69; Cells = @("Elf", "[25", " -3]", 15, "A", "H", "A", "H", "A", "P", "P", "P", "S", "P", "S", "P", "S", "P", "S") }
  @{ Row = 70; Cells = @("Elf", "[25", " -3]", 13, "A", "H", "A", "H", "A", "P", "S", "P", "S", "P", "S", "P", "S") }
  @{ Row = 71; Cells = @("Magician", "[15", " -2]", 15, "S", "A", "A", "A", "P", "A", "S", "H", "A", "H", "P", "A", "A", "P", "S") }
  @{ Row = 72; Cells = @("Magician", "[-4", " 6]", 16, "A", "A", "S", "A", "A", "A", "S", "H", "A", "H", "P", "S", "A", "S", "S", "S") }
  @{ Row = 73; Cells = @("Magician", "[12", " -1]", 13, "A", "A", "A", "A", "A", "H", "S", "H", "A", "P", "A", "P", "S") }
  @{ Row = 74; Cells = @("Magician", "[12", " -2]", 15, "A", "A", "A", "A", "A", "H", "A", "H", "P", "A", "A", "P", "A", "P", "S") }
  @{ Row = 75; Cells = @("Magician", "[0", " 8]", 16, "A", "A", "A", "A", "A", "H", "A", "H", "S", "A", "A", "P", "S", "P", "A", "S") }
  @{ Row = 76; Cells = @("Goblin", "[10", " -2]", 21, "S", "S", "A", "H", "A", "P", "S", "P", "A", "S", "P", "A", "A", "H", "P", "P", "S", "S", "A", "A", "A") }
  @{ Row = 77; Cells = @("Goblin", "[13", " -3]", 11, "P", "A", "A", "A", "P", "A", "S", "S", "A", "S", "A") }
  @{ Row = 78; Cells = @("Goblin", "[10", " -4]", 17, "A", "H", "P", "H", "S", "S", "A", "P", "P", "P", "A", "S", "S", "A", "A", "S", "A") }
  @{ Row = 79; Cells = @("Goblin", "[13", " 0]", 13, "A", "P", "S", "A", "P", "P", "S", "S", "A", "H", "S", "S", "A") }
  @{ Row = 80; Cells = @("Goblin", "[-2", " 2]", 18, "A", "A", "S", "H", "S", "P", "A", "P", "P", "P", "S", "H", "P", "S", "S", "S", "A", "A") }
)

foreach ($r in $newRows) {
  $rowNum = $r.Row
  $cells = $r.Cells
  for ($i = 0; $i -lt $cells.Count; $i++) {
    $ws.Cells.Item($rowNum, $i + 1).Value = $cells[$i]
  }
}